$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    12  = -10.534
    27  = -13.269
    32  = -13.006
    36  = -12.732
    38  = -12.544
    46  = -13.836
    54  = -13.125
    55  = -13.65
    56  = -13.166
    67  = -11.518
    69  = -10.744
    72  = -11.753
    83  = -13.126
    86  = -14.01
    91  = -12.136
    93  = -10.606
    99  = -11.412
    104 = -12.729
}

foreach ($row in $updates.Keys) {
    $ws.Range("C$row").Value = $updates[$row]
}
